# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7067
$ws1.Range("G4").Value = "不可售"
$ws1.Range("F7").Value = 7608
$ws1.Range("F8").Value = 84
$ws1.Range("F9").Value = 203
$ws1.Range("F13").Value = 435
$ws1.Range("F16").Value = 426
$ws1.Range("F20").Value = 5459
$ws1.Range("F21").Value = 140
$ws1.Range("F22").Value = 197
$ws1.Range("F23").Value = 887
$ws1.Range("F24").Value = 226
$ws1.Range("F25").Value = 294

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7067
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F7").Value = 7608
$ws4.Range("F8").Value = 84
$ws4.Range("F9").Value = 203
$ws4.Range("F13").Value = 435
$ws4.Range("F16").Value = 426
$ws4.Range("F21").Value = 5459
$ws4.Range("F23").Value = 140
$ws4.Range("F24").Value = 197
$ws4.Range("F25").Value = 887
$ws4.Range("F26").Value = 226
$ws4.Range("F27").Value = 294
